# "Add files via upload" — re-upload of the workbook with a renamed sheet
# tab and the active-cell selection moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from "推广" to "All features".
$ws.Name = "All features"

# Activate the sheet and move the selection from I20 to C7.
$ws.Activate()
$ws.Range("C7").Select()
